# Arrange UML folders and add sequence diagram paths
#
# The UML "composition" diagram on slide 1 is shifted up and to the left as a
# block (to make room for the new sequence-diagram content below it), the
# shapes are renamed to match PowerPoint's post-move numbering, and the
# leftover empty Title placeholder is removed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Reset the custom "first slide number" back to the PowerPoint default.
$p.PageSetup.FirstSlideNumber = 1

# Target Left/Top (in points) for each diagram shape, in the order they
# appear in the Shapes collection, together with the new shape name that
# PowerPoint assigned when the diagram was rearranged.
$moves = @(
    @{ Index = 1;  Name = "Rectangle 222";           Left = 387.0;              Top = 229.48614173228347 },
    @{ Index = 2;  Name = "Folded Corner 223";        Left = 357.0;              Top = 181.25646219291338 },
    @{ Index = 3;  Name = "TextBox 224";               Left = 357.0;              Top = 181.25653843307086 },
    @{ Index = 4;  Name = "Rectangle 225";            Left = 543.0;              Top = 229.48614173228347 },
    @{ Index = 5;  Name = "Elbow Connector 226";      Left = 477.0;              Top = 244.0267716535433  },
    @{ Index = 6;  Name = "Flowchart: Decision 227";  Left = 459.0;              Top = 235.0267716535433  },
    @{ Index = 7;  Name = "Rectangle 228";            Left = 33.0;               Top = 187.25653843307086 },
    @{ Index = 8;  Name = "Rectangle 229";            Left = 213.0;              Top = 188.1752755905512  },
    @{ Index = 9;  Name = "Elbow Connector 43";       Left = 141.0;              Top = 202.25653843307086 },
    @{ Index = 10; Name = "TextBox 231";               Left = 141.0;              Top = 169.25653843307086 },
    @{ Index = 11; Name = "Flowchart: Decision 232";  Left = 123.0;              Top = 193.25653843307086 },
    @{ Index = 12; Name = "Rectangle 233";            Left = 111.0;              Top = 247.25653843307086 },
    @{ Index = 13; Name = "Elbow Connector 43";       Left = 161.77031716062993; Top = 256.0268503937008  },
    @{ Index = 14; Name = "Flowchart: Decision 235";  Left = 147.0;              Top = 277.2565460330709  },
    @{ Index = 15; Name = "TextBox 236";               Left = 93.0;               Top = 313.2565460330709  }
)

foreach ($m in $moves) {
    $sh = $s.Shapes.Item($m.Index)
    $sh.Name = $m.Name
    $sh.Left = $m.Left
    $sh.Top = $m.Top
}

# The trailing empty "Title 26" placeholder is no longer needed once the
# diagram has been rearranged - remove it.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Title 26") {
        $sh.Delete()
    }
}
